$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("Q1").Value = "test"
$ws.Range("Q1").NumberFormat = "#.00"
